$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # SignIn
$ws2 = $wb.Worksheets.Item(2)  # Client

# --- Update SignIn sheet credentials (demo.com -> example.com, new password) ---
$ws1.Range("B2").Value = "admin@example.com"
$ws1.Range("B3").Value = "client@example.com"
$ws1.Range("C2").Value = "T+wyT5u9cOelDJbBWNgxLw=="
$ws1.Range("C3").Value = "T+wyT5u9cOelDJbBWNgxLw=="

# --- Resize columns B & C on SignIn sheet ---
$ws1.Columns("B").ColumnWidth = 19.428571428571427
$ws1.Columns("C").ColumnWidth = 40.42857142857143

# --- Switch the active/selected sheet from Client to SignIn, update selections ---
$ws2.Activate()
$ws2.Range("C15").Select()
$ws1.Activate()
$ws1.Range("C7").Select()
